# Applies the "excel 2 rows changes" edit:
#  - Row 1 and Row 2 are fully blanked out (contents cleared, styles kept,
#    row height reset to default).
#  - Row 3 is updated to the "Scroll bar missing" bug report (High/Major,
#    status New) and its row height set to 165.
#  - Row 4 is updated to the "HOME PAGE should be enlarged" bug report
#    (Cosmetic, status New) and its row height set to 165.
#  - The data-validation dropdowns on E1:E4 / H1:H4 are removed.
#  - The active selection moves to E3 (no more frozen/scrolled topLeftCell).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 1: clear all content, keep formatting, reset row height ----
$ws.Range("A1:H1").ClearContents() | Out-Null
$ws.Rows(1).AutoFit() | Out-Null

# ---- Row 2: clear all content, keep formatting, reset row height ----
$ws.Range("A2:H2").ClearContents() | Out-Null
$ws.Rows(2).AutoFit() | Out-Null

# ---- Row 3: "Scroll bar missing for navigation panel" bug ----
$ws.Range("E3").Value = "High/Major"
$ws.Range("F3").Value = "Summary:`r`n Scroll bar missing for navigation panel `r`nSteps to Reproduce:`r`nStep 1. Login as Manager`r`n`r`nObservation: `r`nA saparate Scroll bar missing for navigation panel  which leads to empty navigation panel when there are huge no of data in the main screen.`r`n"
$ws.Range("H3").Value = "New"
$ws.Rows(3).RowHeight = 165

# ---- Row 4: "HOME PAGE should be enlarged" bug ----
$ws.Range("E4").Value = "Cosmetic"
$ws.Range("F4").Value = "The word `"HOME PAGE`" should be enlarged.`r`n`r`nSteps to Reproduce:`r`nStep 1. Login as Manager`r`n`r`nObservation: `r`nSince it is the opening page of the application the look and feel of the application improvises by enlargning the size of  `"HOME PAGE`". "
$ws.Range("H4").Value = "New"
$ws.Rows(4).RowHeight = 165

# ---- Remove the list data-validation dropdowns on columns E and H ----
$ws.Range("E1:E4").Validation.Delete()
$ws.Range("H1:H4").Validation.Delete()

# ---- Update the view: selection on E3, no pinned top-left scroll cell ----
$ws.Range("E3").Select() | Out-Null
